$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "2020-09-30 00:00:00"
$ws.Range("I2").Value = 0.12
$ws.Range("K2").Value = 89942100.12
$ws.Range("L2").Value = 5357437.88
$ws.Range("N2").Value = 27.5546035776
$ws.Range("O2").Value = -13.4346335287
$ws.Range("P2").Value = 0.630100861
$ws.Range("Q2").Value = 0.0091204535
$ws.Range("R2").Value = 30.7711064374
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "1"
$ws.Range("AC2").Value = "2020Q3"
$ws.Range("AD2").Value = "2020年 三季报"
$ws.Range("AE2").NumberFormat = "@"
$ws.Range("AE2").Value = "2020"
